# Update the "Marking" row (B11) and "Total" row (B12 / E12) on the
# marksheet so the corrected/total mark figures reflect the new scoring.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking: Right answer value 3 -> 5
$ws.Cells.Item(11, 2).Value = 5

# Total: Right answer total 48 -> 80
$ws.Cells.Item(12, 2).Value = 80

# Total: Correct/total marks label "38/84" -> "80/140"
$ws.Cells.Item(12, 5).Value = "80/140"
